$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.56338
$ws.Range("H2").Value = 1.69014
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5473123333333333
$ws.Range("N2").Value = 1.641937
$ws.Range("O2").Value = 0.0240663820255865
$ws.Range("P2").Value = 0.02406638202558651
$ws.Range("Q2").Value = 0.3083448223533333
$ws.Range("R2").Value = 2.77510340118
$ws.Range("S2").Value = 0.0240663820255865
$ws.Range("T2").Value = 0.02406638202558651

# Row 3
$ws.Range("G3").Value = 0.56338
$ws.Range("H3").Value = 1.69014
$ws.Range("O3").Value = 0.5739535383160298
$ws.Range("P3").Value = 0.5739535383160298
$ws.Range("Q3").Value = 7.353643834913335
$ws.Range("R3").Value = 66.18279451422001
$ws.Range("S3").Value = 0.5739535383160298
$ws.Range("T3").Value = 0.5739535383160298

# Row 4
$ws.Range("G4").Value = 0.56338
$ws.Range("H4").Value = 1.69014
$ws.Range("M4").Value = 9.141742000000001
$ws.Range("N4").Value = 27.425226
$ws.Range("O4").Value = 0.4019800796583838
$ws.Range("P4").Value = 0.4019800796583838
$ws.Range("Q4").Value = 5.15027460796
$ws.Range("R4").Value = 46.35247147164
$ws.Range("S4").Value = 0.4019800796583838
$ws.Range("T4").Value = 0.4019800796583838
